$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Processing Status" sheet lists company names in column A with their
# status in column B. The company previously listed as "Byju" has been
# renamed/corrected to its full legal name: "Arishti CyberTech Private
# Limited". Update the cell in place.
$ws.Range("A2").Value = "Arishti CyberTech Private Limited"
